$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "XTzVQ360"
$ws.Range("B2").Value = 23071736
$ws.Range("C2").Value = "gglvyxp82"
$ws.Range("D2").Value = "vf8#4`$EV"
$ws.Range("F2").Value = "vbGGAtLy"
$ws.Range("G2").Value = "OHsZ"
